# New table and DDL commands Practise
# Adds a new "Sheet4" worksheet (after Sheet3, becoming the active tab) that
# documents the database schema (DBACCOUNT) with six tables: TBLCLIENTS,
# TBLPRODUCTS, TBLEMPLOYEE, TBLSALES, TBLCATEGORIES, TBLCITIES.

$wb = $excel.ActiveWorkbook

# Add the new worksheet at the very end of the workbook (after the current
# last sheet), matching the order tables/columns were originally authored.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)

# Page title
$ws.Range("A1").Value = "DBACCOUNT"

# Table headers (row 3) - products (D) filled before clients (A) to match
# the original authoring order.
$ws.Range("D3").Value = "TBLPRODUCTS"
$ws.Range("A3").Value = "TBLCLIENTS"
$ws.Range("G3").Value = "TBLEMPLOYEE"
$ws.Range("J3").Value = "TBLSALES"
$ws.Range("M3").Value = "TBLCATEGORIES"
$ws.Range("Q3").Value = "TBLCITIES"

# TBLPRODUCTS columns
$ws.Range("D5").Value = "ID"
$ws.Range("D6").Value = "NAME"
$ws.Range("D7").Value = "BRAND"
$ws.Range("D8").Value = "CATEGORYID"
$ws.Range("D9").Value = "COSTPRICE"
$ws.Range("D10").Value = "SELLPRICE"
$ws.Range("D11").Value = "STOCKLEVEL"
$ws.Range("D12").Value = "STATUS"

# TBLEMPLOYEE columns
$ws.Range("G5").Value = "ID"
$ws.Range("G6").Value = "FULLNAME"

# TBLCLIENTS columns
$ws.Range("A5").Value = "ID"
$ws.Range("A6").Value = "FIRSTNAME"
$ws.Range("A7").Value = "LASTNAME"
$ws.Range("A8").Value = "CITYID"
$ws.Range("A9").Value = "BALANCE"

# TBLSALES columns
$ws.Range("J5").Value = "ID"
$ws.Range("J6").Value = "PRODUCTID"
$ws.Range("J7").Value = "CLIENTID"
$ws.Range("J8").Value = "EMPLOYEEID"
$ws.Range("J9").Value = "QUANTITY"
$ws.Range("J10").Value = "TOTAL"
$ws.Range("J11").Value = "DATE"

# TBLCATEGORIES columns
$ws.Range("M5").Value = "ID"
$ws.Range("M6").Value = "TITLE"

# TBLCITIES columns
$ws.Range("Q5").Value = "ID"
$ws.Range("Q6").Value = "LOCATION"

# Match the recorded selection on the new sheet.
$ws.Range("K14").Select()
